# Update forecast values in column B (MSTL) per updated return statistics
# and new scenario generation method.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B2"  = 199.6580963134766
    "B3"  = 198.0741424560547
    "B4"  = 187.6286163330078
    "B5"  = 188.0948181152344
    "B6"  = 186.2721710205078
    "B7"  = 180.1475067138672
    "B8"  = 174.4728851318359
    "B9"  = 177.1927032470703
    "B10" = 188.7744293212891
    "B11" = 204.4593200683594
    "B12" = 203.9915618896484
    "B13" = 226.1586761474609
    "B14" = 226.6010437011719
    "B15" = 222.2899627685547
    "B16" = 190.1591949462891
    "B17" = 175.5337677001953
    "B18" = 164.9937744140625
    "B19" = 150.3785858154297
    "B20" = 136.6065521240234
    "B21" = 125.0766220092773
    "B22" = 127.5506057739258
    "B23" = 113.1779174804688
    "B24" = 114.7959747314453
    "B25" = 102.3314056396484
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
